$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 26-40 (columns A:C) down to rows 27-41 to make room
# for a newly inserted price-list entry at row 26.
for ($r = 40; $r -ge 26; $r--) {
    $dst = $r + 1
    $ws.Cells.Item($dst, 1).Value = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($dst, 2).Value = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($dst, 3).Value = $ws.Cells.Item($r, 3).Value()
}

# New row 26: D82 model entry
$ws.Range("A26").Value = "D82"
$ws.Range("B26").Value = 1025
$ws.Range("C26").Value = 1099

# Update the "Last Update" banner text (new string is added first so it
# lands right before the updated banner string in the shared strings table).
$ws.Range("F4").Value = "Last Update(19-11-2020)"

# Restore the previously selected cell
$ws.Range("K9").Select()
